$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 2041760.1
$ws.Cells.Item(15, 9).Value = 2041760.1
$ws.Cells.Item(15, 11).Value = 6125280.300000001
$ws.Cells.Item(15, 13).Value = -6125111.300000001
$ws.Cells.Item(32, 8).Value = 10696.25
$ws.Cells.Item(32, 10).Value = 7097.25
$ws.Cells.Item(32, 12).Value = 7097.25
$ws.Cells.Item(32, 14).Value = -7749.25
$ws.Cells.Item(43, 8).Value = 10000.5
$ws.Cells.Item(43, 10).Value = 10000
$ws.Cells.Item(43, 12).Value = 10000
$ws.Cells.Item(43, 14).Value = -10138
$ws.Cells.Item(76, 8).Value = 4392.8
$ws.Cells.Item(76, 9).Value = 4226.8237
$ws.Cells.Item(76, 10).Value = 5333.3335
$ws.Cells.Item(76, 11).Value = 4226.8237
$ws.Cells.Item(76, 12).Value = 5333.3335
$ws.Cells.Item(76, 13).Value = -3911.8237
$ws.Cells.Item(76, 14).Value = -5963.3335
$ws.Cells.Item(79, 8).Value = 4392.8
$ws.Cells.Item(79, 9).Value = 4226.8237
$ws.Cells.Item(79, 10).Value = 5333.3335
$ws.Cells.Item(79, 11).Value = 4226.8237
$ws.Cells.Item(79, 12).Value = 5333.3335
$ws.Cells.Item(79, 13).Value = -3134.8237
$ws.Cells.Item(79, 14).Value = -7517.3335
$ws.Cells.Item(80, 8).Value = 13681.091
$ws.Cells.Item(80, 9).Value = 581.6667
$ws.Cells.Item(80, 10).Value = 29400.4
$ws.Cells.Item(80, 11).Value = 1745.0001
$ws.Cells.Item(80, 12).Value = 88201.20000000001
$ws.Cells.Item(80, 13).Value = -747.0001
$ws.Cells.Item(80, 14).Value = -90197.20000000001
$ws.Cells.Item(83, 8).Value = 13681.091
$ws.Cells.Item(83, 9).Value = 581.6667
$ws.Cells.Item(83, 10).Value = 29400.4
$ws.Cells.Item(83, 11).Value = 5235.0003
$ws.Cells.Item(83, 12).Value = 264603.6
$ws.Cells.Item(83, 13).Value = -243.0002999999997
$ws.Cells.Item(83, 14).Value = -274587.6
$ws.Cells.Item(88, 8).Value = 9383
$ws.Cells.Item(88, 10).Value = 13624.5
$ws.Cells.Item(88, 12).Value = 13624.5
$ws.Cells.Item(88, 14).Value = -14436.5
$ws.Cells.Item(91, 8).Value = 9383
$ws.Cells.Item(91, 10).Value = 13624.5
$ws.Cells.Item(91, 12).Value = 13624.5
$ws.Cells.Item(91, 14).Value = -16432.5
$ws.Cells.Item(109, 8).Value = 87342
$ws.Cells.Item(109, 10).Value = 87342
$ws.Cells.Item(109, 12).Value = 87342
$ws.Cells.Item(109, 14).Value = -90116
$ws.Cells.Item(115, 8).Value = 1081.4286
$ws.Cells.Item(115, 9).Value = 1214
$ws.Cells.Item(115, 10).Value = 750
$ws.Cells.Item(115, 11).Value = 3642
$ws.Cells.Item(115, 12).Value = 2250
$ws.Cells.Item(115, 13).Value = -2075
$ws.Cells.Item(115, 14).Value = -5384
$ws.Cells.Item(137, 8).Value = 13747.793
$ws.Cells.Item(137, 9).Value = 15617.174
$ws.Cells.Item(137, 11).Value = 46851.522
$ws.Cells.Item(137, 13).Value = -44301.522
$ws.Cells.Item(141, 8).Value = 851.2941
$ws.Cells.Item(141, 10).Value = 323.5
$ws.Cells.Item(141, 12).Value = 970.5
$ws.Cells.Item(141, 14).Value = -11330.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 30189.361
$ws.Cells.Item(32, 9).Value = 22703.463
$ws.Cells.Item(32, 11).Value = 22703.463
$ws.Cells.Item(32, 13).Value = -22416.463
$ws.Cells.Item(45, 8).Value = 4771.6113
$ws.Cells.Item(45, 9).Value = 3199.625
$ws.Cells.Item(45, 10).Value = 6029.2
$ws.Cells.Item(45, 11).Value = 3199.625
$ws.Cells.Item(45, 12).Value = 6029.2
$ws.Cells.Item(45, 13).Value = -2822.625
$ws.Cells.Item(45, 14).Value = -6783.2
$ws.Cells.Item(61, 8).Value = 4207.2563
$ws.Cells.Item(61, 9).Value = 1249.4517
$ws.Cells.Item(61, 11).Value = 1249.4517
$ws.Cells.Item(61, 13).Value = -1037.4517
$ws.Cells.Item(63, 8).Value = 2071.625
$ws.Cells.Item(63, 10).Value = 3666.3333
$ws.Cells.Item(63, 12).Value = 3666.3333
$ws.Cells.Item(63, 14).Value = -5038.3333
$ws.Cells.Item(66, 8).Value = 2071.625
$ws.Cells.Item(66, 10).Value = 3666.3333
$ws.Cells.Item(66, 12).Value = 18331.6665
$ws.Cells.Item(66, 14).Value = -25195.6665
$ws.Cells.Item(74, 8).Value = 125218.695
$ws.Cells.Item(74, 9).Value = 147056.81
$ws.Cells.Item(74, 11).Value = 147056.81
$ws.Cells.Item(74, 13).Value = -146182.81
$ws.Cells.Item(77, 8).Value = 125218.695
$ws.Cells.Item(77, 9).Value = 147056.81
$ws.Cells.Item(77, 11).Value = 735284.05
$ws.Cells.Item(77, 13).Value = -730916.05
$ws.Cells.Item(102, 8).Value = 5341.0435
$ws.Cells.Item(102, 9).Value = 4229.25
$ws.Cells.Item(102, 10).Value = 7882.2856
$ws.Cells.Item(102, 11).Value = 4229.25
$ws.Cells.Item(102, 12).Value = 7882.2856
$ws.Cells.Item(102, 13).Value = -2607.25
$ws.Cells.Item(102, 14).Value = -11126.2856
$ws.Cells.Item(136, 8).Value = 4207.2563
$ws.Cells.Item(136, 9).Value = 1249.4517
$ws.Cells.Item(136, 11).Value = 3748.3551
$ws.Cells.Item(136, 13).Value = -1198.3551

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1765.5807
$ws.Cells.Item(58, 9).Value = 1684.375
$ws.Cells.Item(58, 10).Value = 2044
$ws.Cells.Item(58, 11).Value = 1684.375
$ws.Cells.Item(58, 12).Value = 2044
$ws.Cells.Item(58, 13).Value = -1481.375
$ws.Cells.Item(58, 14).Value = -2450
$ws.Cells.Item(62, 8).Value = 8476.1
$ws.Cells.Item(62, 10).Value = 8356.714
$ws.Cells.Item(62, 12).Value = 8356.714
$ws.Cells.Item(62, 14).Value = -9604.714
$ws.Cells.Item(65, 8).Value = 8476.1
$ws.Cells.Item(65, 10).Value = 8356.714
$ws.Cells.Item(65, 12).Value = 41783.57
$ws.Cells.Item(65, 14).Value = -48023.57
$ws.Cells.Item(70, 8).Value = 41625
$ws.Cells.Item(70, 10).Value = 41625
$ws.Cells.Item(70, 12).Value = 41625
$ws.Cells.Item(70, 14).Value = -42255
$ws.Cells.Item(73, 8).Value = 41625
$ws.Cells.Item(73, 10).Value = 41625
$ws.Cells.Item(73, 12).Value = 41625
$ws.Cells.Item(73, 14).Value = -43809
$ws.Cells.Item(122, 8).Value = 2575.1875
$ws.Cells.Item(122, 9).Value = 2465.12
$ws.Cells.Item(122, 10).Value = 2968.2856
$ws.Cells.Item(122, 11).Value = 7395.36
$ws.Cells.Item(122, 12).Value = 8904.856800000001
$ws.Cells.Item(122, 13).Value = -4945.36
$ws.Cells.Item(122, 14).Value = -13804.8568
$ws.Cells.Item(132, 8).Value = 28186.514
$ws.Cells.Item(132, 9).Value = 32206.312
$ws.Cells.Item(132, 10).Value = 2459.8
$ws.Cells.Item(132, 11).Value = 96618.936
$ws.Cells.Item(132, 12).Value = 7379.400000000001
$ws.Cells.Item(132, 13).Value = -94088.936
$ws.Cells.Item(132, 14).Value = -12439.4
$ws.Cells.Item(136, 8).Value = 1765.5807
$ws.Cells.Item(136, 9).Value = 1684.375
$ws.Cells.Item(136, 10).Value = 2044
$ws.Cells.Item(136, 11).Value = 5053.125
$ws.Cells.Item(136, 12).Value = 6132
$ws.Cells.Item(136, 13).Value = -2503.125
$ws.Cells.Item(136, 14).Value = -11232

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(70, 8).Value = 4695.6763
$ws.Cells.Item(70, 9).Value = 3487.6
$ws.Cells.Item(70, 10).Value = 4903.9653
$ws.Cells.Item(70, 11).Value = 10462.8
$ws.Cells.Item(70, 12).Value = 14711.8959
$ws.Cells.Item(70, 13).Value = -10147.8
$ws.Cells.Item(70, 14).Value = -15341.8959
$ws.Cells.Item(73, 8).Value = 4695.6763
$ws.Cells.Item(73, 9).Value = 3487.6
$ws.Cells.Item(73, 10).Value = 4903.9653
$ws.Cells.Item(73, 11).Value = 10462.8
$ws.Cells.Item(73, 12).Value = 14711.8959
$ws.Cells.Item(73, 13).Value = -9370.799999999999
$ws.Cells.Item(73, 14).Value = -16895.8959

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 8013.1904
$ws.Cells.Item(80, 9).Value = 4939.8887
$ws.Cells.Item(80, 10).Value = 10318.167
$ws.Cells.Item(80, 11).Value = 4939.8887
$ws.Cells.Item(80, 12).Value = 10318.167
$ws.Cells.Item(80, 13).Value = -3941.8887
$ws.Cells.Item(80, 14).Value = -12314.167
$ws.Cells.Item(83, 8).Value = 8013.1904
$ws.Cells.Item(83, 9).Value = 4939.8887
$ws.Cells.Item(83, 10).Value = 10318.167
$ws.Cells.Item(83, 11).Value = 24699.4435
$ws.Cells.Item(83, 12).Value = 51590.835
$ws.Cells.Item(83, 13).Value = -19707.4435
$ws.Cells.Item(83, 14).Value = -61574.835
$ws.Cells.Item(97, 8).Value = 3129
$ws.Cells.Item(97, 10).Value = 3213.9092
$ws.Cells.Item(97, 12).Value = 3213.9092
$ws.Cells.Item(97, 14).Value = -4205.9092
$ws.Cells.Item(102, 8).Value = 7835.8696
$ws.Cells.Item(102, 9).Value = 9120.294
$ws.Cells.Item(102, 11).Value = 9120.294
$ws.Cells.Item(102, 13).Value = -7498.294
$ws.Cells.Item(126, 8).Value = 3930.25
$ws.Cells.Item(126, 9).Value = 2499.5
$ws.Cells.Item(126, 10).Value = 4645.625
$ws.Cells.Item(126, 11).Value = 7498.5
$ws.Cells.Item(126, 12).Value = 13936.875
$ws.Cells.Item(126, 13).Value = -5028.5
$ws.Cells.Item(126, 14).Value = -18876.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 3347.2354
$ws.Cells.Item(68, 10).Value = 4250.6665
$ws.Cells.Item(68, 12).Value = 4250.6665
$ws.Cells.Item(68, 14).Value = -5748.6665
$ws.Cells.Item(71, 8).Value = 3347.2354
$ws.Cells.Item(71, 10).Value = 4250.6665
$ws.Cells.Item(71, 12).Value = 21253.3325
$ws.Cells.Item(71, 14).Value = -28741.3325
$ws.Cells.Item(93, 8).Value = 1789.1666
$ws.Cells.Item(93, 9).Value = 1480.75
$ws.Cells.Item(93, 11).Value = 1480.75
$ws.Cells.Item(93, 13).Value = -232.75
$ws.Cells.Item(132, 8).Value = 1632.6349
$ws.Cells.Item(132, 9).Value = 1181.5818
$ws.Cells.Item(132, 11).Value = 3544.7454
$ws.Cells.Item(132, 13).Value = -1014.7454

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 1011739.3
$ws.Cells.Item(5, 10).Value = 1011739.3
$ws.Cells.Item(5, 12).Value = 1011739.3
$ws.Cells.Item(5, 14).Value = -1011963.3
$ws.Cells.Item(62, 8).Value = 3715
$ws.Cells.Item(62, 9).Value = 3358.2
$ws.Cells.Item(62, 10).Value = 5499
$ws.Cells.Item(62, 11).Value = 3358.2
$ws.Cells.Item(62, 12).Value = 5499
$ws.Cells.Item(62, 13).Value = -2734.2
$ws.Cells.Item(62, 14).Value = -6747
$ws.Cells.Item(65, 8).Value = 3715
$ws.Cells.Item(65, 9).Value = 3358.2
$ws.Cells.Item(65, 10).Value = 5499
$ws.Cells.Item(65, 11).Value = 16791
$ws.Cells.Item(65, 12).Value = 27495
$ws.Cells.Item(65, 13).Value = -13671
$ws.Cells.Item(65, 14).Value = -33735
$ws.Cells.Item(122, 8).Value = 10445758
$ws.Cells.Item(122, 9).Value = 12228410
$ws.Cells.Item(122, 10).Value = 4507.857
$ws.Cells.Item(122, 11).Value = 36685230
$ws.Cells.Item(122, 12).Value = 13523.571
$ws.Cells.Item(122, 13).Value = -36682780
$ws.Cells.Item(122, 14).Value = -18423.571
$ws.Cells.Item(126, 8).Value = 3550.5715
$ws.Cells.Item(126, 9).Value = 2071.3
$ws.Cells.Item(126, 11).Value = 6213.900000000001
$ws.Cells.Item(126, 13).Value = -3743.900000000001
$ws.Cells.Item(132, 8).Value = 6599382.5
$ws.Cells.Item(132, 9).Value = 10029550
$ws.Cells.Item(132, 10).Value = 2906.3076
$ws.Cells.Item(132, 11).Value = 30088650
$ws.Cells.Item(132, 12).Value = 8718.9228
$ws.Cells.Item(132, 13).Value = -30086120
$ws.Cells.Item(132, 14).Value = -13778.9228
